$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (e.g. "30.321.41", "1.001") that Excel would
# otherwise auto-convert to numbers. Force the Price column to Text format
# first so the new values are stored as text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.315.10'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '1.932.17'
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '251.43'
$ws.Range("E5").Value = '  +1.40%  '
$ws.Range("D6").Value = '0.7143'
$ws.Range("E6").Value = '  -0.84%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '0.3264'
$ws.Range("E8").Value = '  -0.75%  '
$ws.Range("D9").Value = '27.38'
$ws.Range("E9").Value = '  +2.59%  '
$ws.Range("D10").Value = '0.07196'
$ws.Range("E10").Value = '  +5.21%  '
$ws.Range("D11").Value = '0.7992'
$ws.Range("E11").Value = '  -1.23%  '
$ws.Range("D12").Value = '0.08096'
$ws.Range("E12").Value = '  +1.46%  '
$ws.Range("D13").Value = '1.929.99'
$ws.Range("E13").Value = '  -0.22%  '
$ws.Range("D14").Value = '5.435'
$ws.Range("E14").Value = '  -0.28%  '
$ws.Range("D15").Value = '94.82'
$ws.Range("E15").Value = '  -0.13%  '
$ws.Range("D16").Value = '14.84'
$ws.Range("E16").Value = '  +1.53%  '
$ws.Range("D17").Value = '30.328.21'
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("D18").Value = '253.34'
$ws.Range("E18").Value = '  -3.94%  '
$ws.Range("D19").Value = '0.000008112'
$ws.Range("E19").Value = '  +1.87%  '
$ws.Range("D20").Value = '5.804'
$ws.Range("E20").Value = '  -0.72%  '
$ws.Range("D21").Value = '2.182.76'
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").Value = '6.920'
$ws.Range("E24").Value = '  -0.23%  '
$ws.Range("D25").Value = '9.708'
$ws.Range("E25").Value = '  -0.59%  '
$ws.Range("D26").Value = '164.88'
$ws.Range("D27").Value = '19.27'
$ws.Range("E27").Value = '  +0.96%  '
$ws.Range("D28").Value = '2.316'
$ws.Range("E28").Value = '  -1.44%  '
$ws.Range("D29").Value = '0.1282'
$ws.Range("E29").Value = '  -4.79%  '
$ws.Range("D30").Value = '1.360'
$ws.Range("E30").Value = '  -0.55%  '
$ws.Range("D31").Value = '1.544'
$ws.Range("E31").Value = '  -0.83%  '
$ws.Range("D32").Value = '4.433'
$ws.Range("E32").Value = '  +0.18%  '
$ws.Range("D33").Value = '4.206'
$ws.Range("E33").Value = '  -0.72%  '
$ws.Range("D34").Value = '0.05215'
$ws.Range("E34").Value = '  +2.24%  '
$ws.Range("D35").Value = '1.268'
$ws.Range("E35").Value = '  +4.67%  '
$ws.Range("D36").Value = '0.7500'
$ws.Range("E36").Value = '  +0.16%  '
$ws.Range("D37").Value = '2.767'
$ws.Range("E37").Value = '  +0.99%  '
$ws.Range("D38").Value = '0.01966'
$ws.Range("E38").Value = '  +0.83%  '
$ws.Range("D39").Value = '2.801'
$ws.Range("E39").Value = '  -0.85%  '
$ws.Range("D40").Value = '79.07'
$ws.Range("E40").Value = '  -2.84%  '
$ws.Range("D41").Value = '6.440'
$ws.Range("E41").Value = '  -2.32%  '
$ws.Range("D42").Value = '0.4535'
$ws.Range("E42").Value = '  +0.83%  '
$ws.Range("D43").Value = '2.027'
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("D44").Value = '0.8422'
$ws.Range("E44").Value = '  +0.55%  '
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("D46").Value = '102.03'
$ws.Range("E46").Value = '  -0.50%  '
$ws.Range("D47").Value = '9.816'
$ws.Range("E47").Value = '  +0.46%  '
$ws.Range("D48").Value = '7.452'
$ws.Range("E48").Value = '  +1.16%  '
$ws.Range("D49").Value = '36.72'
$ws.Range("E49").Value = '  +0.89%  '
$ws.Range("D50").Value = '0.06095'
$ws.Range("E50").Value = '  +2.53%  '
$ws.Range("D51").Value = '0.4185'
$ws.Range("E51").Value = '  +1.22%  '
